# Data_Login.xlsx — "Add files via upload" edit
# Adds 7 more login rows (Test_URL / User_Name / Password) below the
# existing sample row, turns each new User_Name into a mailto: hyperlink,
# widens column B to fit the longer addresses, and moves the active
# selection onto the newly added Password column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$testUrl  = $ws.Range("A2").Value()
$password = $ws.Range("C2").Value()

$emails = @(
    "barakathk@hlwe.com",
    "leems@hlwe.com",
    "swaran@hlwe.com",
    "francis@hlwe.com",
    "linda@hlwe.com",
    "adele@hlwe.com",
    "tancs@hlwe.com"
)

$row = 3
foreach ($email in $emails) {
    $ws.Cells.Item($row, 1).Value = $testUrl
    $ws.Cells.Item($row, 2).Value = $email
    $ws.Cells.Item($row, 3).Value = $password

    $bCell = $ws.Cells.Item($row, 2)
    $ws.Hyperlinks.Add($bCell, "mailto:" + $email, [Type]::Missing, [Type]::Missing, "mailto:" + $email)
    # TextToDisplay above overwrote the cell with "mailto:<email>" — put the
    # plain address back as the visible cell text (hyperlink stays attached).
    $bCell.Value = $email
    # The hyperlink auto-applies the built-in "Hyperlink" style (blue/underline);
    # strip that back off so the cell keeps its original, unstyled look.
    $bCell.ClearFormats()

    $row = $row + 1
}

# Column B now holds longer addresses than before — resize to fit.
$ws.Columns("B").AutoFit()

# Leave the selection on the freshly added Password entries, like the author did.
$ws.Range("C7:C9").Select() | Out-Null

Write-Output "done"
